# Changed SummaryStats so that the last value has to be a detected result.
# Updates the "Last Cr" (X) and "Last Date" (Y) columns on the
# "Alluvial for Mapping" worksheet for rows 11-17 so that the last
# reported chromium result reflects a detected (non "No Detect") sample,
# mirroring the corresponding Max Cr / Max Date values (columns V/W) for
# most wells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alluvial for Mapping")

# These columns store plain text (e.g. "5.76", "2010-07-06"), not real
# numbers/dates, matching the rest of the sheet (and columns V/W in
# particular). Mark the cells as Text before assigning so Excel keeps the
# values as literal strings instead of auto-converting them to a number or
# a date serial value.
$updates = @(
    @{ Row = 11; X = "5.76";            Y = "2010-07-06" },
    @{ Row = 12; X = "5.44";            Y = "2010-07-07" },
    @{ Row = 13; X = "3.43";            Y = "2018-06-18" },
    @{ Row = 14; X = "0.692";           Y = "2002-05-22" },
    @{ Row = 15; X = "No Detect Data";  Y = "No Detect Data" },
    @{ Row = 16; X = "7.76";            Y = "2010-07-26" },
    @{ Row = 17; X = "3.01";            Y = "2010-06-08" }
)

foreach ($u in $updates) {
    $xCell = $ws.Range("X$($u.Row)")
    $yCell = $ws.Range("Y$($u.Row)")

    # "No Detect Data" is unambiguous text already; the numeric/date-like
    # values need an explicit text format so they aren't reinterpreted.
    if ($u.X -notmatch "^[A-Za-z]") {
        $xCell.NumberFormat = "@"
    }
    if ($u.Y -notmatch "^[A-Za-z]") {
        $yCell.NumberFormat = "@"
    }

    $xCell.Value = $u.X
    $yCell.Value = $u.Y
}
